# Refresh cryptos list: prices (col D) and 1h volume deltas (col E) updated per source diff.
# Row 45/46 also swap coin order (FraxShare <-> EnergySwap) with new data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "`'39.948.26"
$ws.Range("E2").Value = "`'  +0.61%  "
$ws.Range("D3").Value = "`'2.219.92"
$ws.Range("E3").Value = "`'  +0.35%  "
$ws.Range("D4").Value = "`'1.00"
$ws.Range("E4").Value = "`'  -0.04%  "
$ws.Range("D5").Value = "`'292.43"
$ws.Range("E5").Value = "`'  -0.01%  "
$ws.Range("D6").Value = "`'87.96"
$ws.Range("E6").Value = "`'  +2.25%  "
$ws.Range("E7").Value = "`'  +0.05%  "
$ws.Range("E8").Value = "`'  -0.06%  "
$ws.Range("E9").Value = "`'  -0.44%  "
$ws.Range("D10").Value = "`'30.68"
$ws.Range("E10").Value = "`'  +0.19%  "
$ws.Range("D11").Value = "`'0.0782"
$ws.Range("E11").Value = "`'  -0.30%  "
$ws.Range("D12").Value = "`'50.02"
$ws.Range("E12").Value = "`'  +5.27%  "
$ws.Range("E13").Value = "`'  +2.46%  "
$ws.Range("D14").Value = "`'6.47"
$ws.Range("E14").Value = "`'  +2.16%  "
$ws.Range("D15").Value = "`'2.561.74"
$ws.Range("E15").Value = "`'  +0.20%  "
$ws.Range("D16").Value = "`'13.80"
$ws.Range("E16").Value = "`'  -1.44%  "
$ws.Range("D17").Value = "`'2.195.77"
$ws.Range("E17").Value = "`'  -0.82%  "
$ws.Range("E18").Value = "`'  +0.39%  "
$ws.Range("D19").Value = "`'39.901.05"
$ws.Range("E19").Value = "`'  +0.54%  "
$ws.Range("D20").Value = "`'0.0₃0888"
$ws.Range("E20").Value = "`'  +0.92%  "
$ws.Range("D21").Value = "`'11.16"
$ws.Range("E21").Value = "`'  -0.55%  "
$ws.Range("E22").Value = "`'  -0.90%  "
$ws.Range("D23").Value = "`'65.63"
$ws.Range("E23").Value = "`'  +0.26%  "
$ws.Range("D24").Value = "`'237.28"
$ws.Range("E24").Value = "`'  +0.87%  "
$ws.Range("E26").Value = "`'  +0.07%  "
$ws.Range("E27").Value = "`'  +0.11%  "
$ws.Range("D28").Value = "`'23.21"
$ws.Range("E28").Value = "`'  +2.16%  "
$ws.Range("E29").Value = "`'  +0.45%  "
$ws.Range("E30").Value = "`'  -6.74%  "
$ws.Range("D31").Value = "`'157.04"
$ws.Range("E31").Value = "`'  +3.65%  "
$ws.Range("D32").Value = "`'31.94"
$ws.Range("E32").Value = "`'  -2.53%  "
$ws.Range("D33").Value = "`'0.999"
$ws.Range("E33").Value = "`'  +0.07%  "
$ws.Range("D34").Value = "`'4.98"
$ws.Range("E34").Value = "`'  +1.11%  "
$ws.Range("D35").Value = "`'0.0713"
$ws.Range("E35").Value = "`'  -0.45%  "
$ws.Range("D36").Value = "`'2.96"
$ws.Range("E36").Value = "`'  +6.14%  "
$ws.Range("D37").Value = "`'2.34"
$ws.Range("E37").Value = "`'  -1.38%  "
$ws.Range("E38").Value = "`'  -0.18%  "
$ws.Range("D39").Value = "`'0.0989"
$ws.Range("E39").Value = "`'  +0.26%  "
$ws.Range("E40").Value = "`'  +1.95%  "
$ws.Range("D41").Value = "`'15.37"
$ws.Range("E41").Value = "`'  -2.99%  "
$ws.Range("D42").Value = "`'2.108.40"
$ws.Range("E42").Value = "`'  +2.00%  "
$ws.Range("E43").Value = "`'  -0.78%  "
$ws.Range("E44").Value = "`'  +0.97%  "
$ws.Range("B45").Value = "`'EnergySwap"
$ws.Range("C45").Value = "`'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "`'17.91"
$ws.Range("E45").Value = "`'  +1.48%  "
$ws.Range("B46").Value = "`'FraxShare"
$ws.Range("C46").Value = "`'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "`'9.96"
$ws.Range("E46").Value = "`'  +0.09%  "
$ws.Range("E47").Value = "`'  -0.41%  "
$ws.Range("E48").Value = "`'  +3.76%  "
$ws.Range("D49").Value = "`'2.429.77"
$ws.Range("E49").Value = "`'  -0.14%  "
$ws.Range("E50").Value = "`'  +2.88%  "
$ws.Range("E51").Value = "`'  -0.24%  "
